# Quarterly indexing esoteric bug-fix operation
#
# Column A holds a date stamp per row that is supposed to mark the
# "as-of" reporting date for each quarterly forecast row, but it was
# being indexed off the *start* of the quarter's first month instead of
# correctly landing mid-way through the month that follows it. Re-stamp
# every populated row in column A (the header row, row 1, is untouched)
# by pushing the stored date forward to the 15th of the following month.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $d = $cell.Value()
    if ($null -eq $d) { continue }

    $monthStart = Get-Date -Year $d.Year -Month $d.Month -Day 1 -Hour 0 -Minute 0 -Second 0
    $nextMonth = $monthStart.AddMonths(1)
    $fixed = Get-Date -Year $nextMonth.Year -Month $nextMonth.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value = $fixed
}
